# Add a new "Sheet2" worksheet at the end of the workbook, populate it with
# test data (plain table, an Excel Table/ListObject, and a block used by a
# defined name), then register the defined name.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet2"

# --- simple data block: A1:C6 ----------------------------------------------
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "e"
$ws.Range("C2").Value = 1.1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "d"
$ws.Range("C3").Value = 2.2

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "c"
$ws.Range("C4").Value = 3.3

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "b"
$ws.Range("C5").Value = 3.145612

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "a"
$ws.Range("C6").Value = 0.98

# --- table block: G11:I17 (becomes Table1) ----------------------------------
$ws.Range("G11").Value = "tabc1"
$ws.Range("H11").Value = "tabc2"
$ws.Range("I11").Value = "tabc3"

$ws.Range("G12").Value = "a1"
$ws.Range("H12").Value = $true
$ws.Range("I12").Value = 1.1

$ws.Range("G13").Value = "a2"
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = 2.2

$ws.Range("G14").Value = "a3"
$ws.Range("H14").Value = $true
$ws.Range("I14").Value = 3.3

$ws.Range("G15").Value = "a4"
$ws.Range("H15").Value = $true
$ws.Range("I15").Value = 4.4

$ws.Range("G16").Value = "a5"
$ws.Range("H16").Value = $false
$ws.Range("I16").Value = 5.56

$ws.Range("G17").Value = "a6"
$ws.Range("H17").Value = $true
$ws.Range("I17").Value = 0.9999

# turn G11:I17 into an actual Excel Table named "Table1"
$tbl = $ws.ListObjects.Add(1, $ws.Range("G11:I17"), $null, 1)
$tbl.Name = "Table1"

# --- named-range block: C20:E24 ---------------------------------------------
$ws.Range("C20").Value = "alpha"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 5

$ws.Range("C21").Value = "beta"
$ws.Range("D21").Value = 2.2
$ws.Range("E21").Value = 6

$ws.Range("C22").Value = "charlie"
$ws.Range("D22").Value = 3.3
$ws.Range("E22").Value = 7

$ws.Range("C23").Value = "delta"
$ws.Range("D23").Value = 4.4
$ws.Range("E23").Value = 8

$ws.Range("C24").Value = "echo"
$ws.Range("D24").Value = 5.5
$ws.Range("E24").Value = 9

# register the workbook-level defined name used by the data above
$wb.Names.Add("TestNamedRange", "=Sheet2!`$C`$20:`$E`$24")

# leave the selection where the original file had it
$ws.Range("C20").Select()
